# The post at row 390 (NASA / James Webb Space Telescope tweet) was removed.
# Deleting the entire row shifts all subsequent rows up by one, which matches
# the target diff (rows 391-583 become 390-582, and the sheet dimension
# shrinks from A1:C583 to A1:C582).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("390:390").Delete() | Out-Null
